# The commit swaps the contents of ppt/theme/theme1.xml ("Integral" theme,
# used by the deck's single Slide Master) and ppt/theme/theme2.xml ("Office
# Theme", only ever linked from the Notes Master). Font scheme and format
# scheme (fills/lines/effects) are byte-identical between the two theme
# parts, so the only substantive difference is the 12-colour theme colour
# scheme (and the cosmetic theme/colour-scheme "name" attributes, which
# PowerPoint's object model does not expose as writable).
#
# The PowerPoint COM object model only ever surfaces a single addressable
# theme/colour-scheme for this deck (Master.Theme / NotesMaster.Theme /
# HandoutMaster.Theme all resolve to the same underlying theme part, because
# there is only one Slide Master) so theme2.xml itself cannot be reached
# through the object model. We therefore apply the reachable, concrete part
# of the edit: re-point the live theme's ThemeColorScheme entries at the
# "Office Theme" palette that theme2.xml carried before the swap.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

# Index -> (scheme slot, old "Integral" value, new "Office Theme" value)
#  1 dk1      000000 -> 000000 (unchanged)
#  2 lt1      FFFFFF -> FFFFFF (unchanged)
#  3 dk2      455F51 -> 44546A
#  4 lt2      E3DED1 -> E7E6E6
#  5 accent1  99CB38 -> 5B9BD5
#  6 accent2  63A537 -> ED7D31
#  7 accent3  E6D024 -> A5A5A5
#  8 accent4  CC9700 -> FFC000
#  9 accent5  4EB3CF -> 4472C4
# 10 accent6  378DA6 -> 70AD47
# 11 hlink    6B9F25 -> 0563C1
# 12 folHlink B26B02 -> 954F72
$colors.Item(1).RGB  = 0
$colors.Item(2).RGB  = 16777215
$colors.Item(3).RGB  = 6968388
$colors.Item(4).RGB  = 15132391
$colors.Item(5).RGB  = 13998939
$colors.Item(6).RGB  = 3243501
$colors.Item(7).RGB  = 10855845
$colors.Item(8).RGB  = 49407
$colors.Item(9).RGB  = 12874308
$colors.Item(10).RGB = 4697456
$colors.Item(11).RGB = 12673797
$colors.Item(12).RGB = 7491477

# Best-effort rename (PowerPoint's object model treats these as read-only in
# practice, but set them in case the host honours it).
try { $theme.Name = "Office Theme" } catch {}
try { $colors.Name = "Office" } catch {}
